$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "26.887.72"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.642.27"
$ws.Range("E3").Value = "  -0.15%  "
Set-TextValue "D4" "1.01"
$ws.Range("E4").Value = "  -0.38%  "
Set-TextValue "D5" "219.29"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.21%  "
Set-TextValue "D9" "0.0624"
$ws.Range("E9").Value = "  -0.51%  "
Set-TextValue "D10" "19.27"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.871.31"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "1.640.72"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -0.43%  "
Set-TextValue "D15" "0.528"
$ws.Range("E15").Value = "  +0.09%  "
Set-TextValue "D16" "65.60"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "26.865.77"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  -0.57%  "
Set-TextValue "D19" "216.52"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  +5.20%  "
Set-TextValue "D23" "2.42"
$ws.Range("E23").Value = "  -1.43%  "
Set-TextValue "D24" "9.21"
$ws.Range("E24").Value = "  -1.39%  "
Set-TextValue "D25" "148.09"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +1.50%  "
Set-TextValue "D29" "15.81"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  -0.20%  "
Set-TextValue "D34" "1.56"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").Value = "1.271.32"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -2.07%  "
Set-TextValue "D38" "0.532"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  -0.30%  "
Set-TextValue "D41" "0.807"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "1.781.96"
$ws.Range("E43").Value = "  -0.40%  "
Set-TextValue "D44" "92.63"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D45" "2.08"
$ws.Range("E45").Value = "  -7.04%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "60.96"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("E49").Value = "  -0.77%  "
Set-TextValue "D50" "0.0969"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").Value = "  -0.30%  "
